# Regenerate s_vals data (B,C,D,E,G columns) to reflect filtering of save games.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, TB(B), d2S(C), K(D), IP(E), sum(G)
$data = @(
    @(2, 0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.094976487407548),
    @(3, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(4, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(5, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(6, 0.127881588408715, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.449980674824537),
    @(7, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(8, 3.230985683306322, 10.29869402782916, 3.900430680208489, 8.660232485948974, 26.09034287729295),
    @(9, 0.3048080303191223, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1.001517020209437),
    @(10, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(11, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(12, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(13, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(14, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(15, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(16, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(17, 0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082),
    @(18, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(19, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(20, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 13.71653804550039),
    @(21, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(22, 0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.626907116734944),
    @(23, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742),
    @(24, 0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 8.660232485948974, 9.805878329971296),
    @(25, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(26, 0.04763786555579896, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.014732764554632),
    @(27, 0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.35982162114495),
    @(28, 0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.094976487407548),
    @(29, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742),
    @(30, 1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(31, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(32, 1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 11.945164432584),
    @(33, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(34, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(35, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(36, 0.04763786555579896, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.014732764554632),
    @(37, 1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(38, 0.3048080303191223, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1.001517020209437),
    @(39, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(40, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(41, 1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(42, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(43, 0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044),
    @(44, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(45, 1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(46, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(47, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}
